$wb = $excel.ActiveWorkbook

# --- Scenarios sheet: disable ("N") all scenarios except the two LoginTest rows ---
$scenarios = $wb.Worksheets.Item("Scenarios")
for ($r = 4; $r -le 12; $r++) {
    $scenarios.Cells.Item($r, 1).Value = "N"
}
$scenarios.Range("A3:A12").Select()

# --- Parameters sheet: point the run at remote Saucelabs/Appium (iOS Safari) ---
$parameters = $wb.Worksheets.Item("Parameters")
$parameters.Cells.Item(2, 2).Value = "remote"

$parameters.Cells.Item(14, 1).Value = "appiumVersion"
$parameters.Cells.Item(15, 1).Value = "deviceName"
$parameters.Cells.Item(16, 1).Value = "deviceOrientation"
$parameters.Cells.Item(17, 1).Value = "platformVersion"
$parameters.Cells.Item(18, 1).Value = "platformName"

$parameters.Cells.Item(14, 2).Value = "1.6.4"
$parameters.Cells.Item(15, 2).Value = "iPhone 7 Simulator"
$parameters.Cells.Item(16, 2).Value = "portrait"
$parameters.Cells.Item(17, 2).Value = "10.3"
$parameters.Cells.Item(18, 2).Value = "iOS"

$parameters.Cells.Item(4, 2).Value = "iphone"
$parameters.Cells.Item(3, 2).Value = "Safari"

$parameters.Range("B4").Select()

# --- Configs sheet: no content change, just cursor/selection moves ---
$configs = $wb.Worksheets.Item("Configs")
$configs.Activate()
$configs.Range("E26").Select()

# Return focus to the Parameters tab (it stays the active sheet in the workbook)
$parameters.Activate()
$parameters.Range("B4").Select()
